$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1; everything below (including data rows and
# the two hyperlinked cells) shifts down by one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row with plain (unstyled) reference info.
$ws.Cells.Item(1, 1).Value = "Reference"
$ws.Cells.Item(1, 2).Value = "UNAM"

# The row insert does not relocate the existing hyperlink anchors, so
# recreate them at their shifted positions (old B17/B19 -> new B18/B20).
$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("B18"), "https://en.wikipedia.org/wiki/Michoac%C3%A1n", "", "Michoacán", "https://en.wikipedia.org/wiki/Michoac%C3%A1n")
$null = $ws.Hyperlinks.Add($ws.Range("B20"), "https://en.wikipedia.org/wiki/Veracruz", "", "Veracruz", "https://en.wikipedia.org/wiki/Veracruz")

# Match the target selection state.
$ws.Range("B1").Select()
